# Commit message: "Accepted all changes for 2.23.0 release notes"
#
# The document contained tracked changes (insertions and one deletion,
# all authored by "Menon, Sunita (NIH/NCI) [C]") describing the
# HPCDATAMGM-1629 release-note bullet. This change simply accepts every
# tracked revision in the document, which:
#   - keeps the inserted text/runs (dropping the <w:ins> wrappers)
#   - removes the text covered by the <w:del> markup
#   - resolves the <w:rPrChange>/formatting-change tracking, leaving
#     just the final formatting
#   - merges the paragraph whose mark was itself a tracked insertion
#     into the following paragraph (standard Word behavior when an
#     inserted paragraph mark is accepted)
#
# This mirrors choosing Review > Accept > Accept All Changes in Word.

$d = $word.ActiveDocument
$d.AcceptAllRevisions()
